$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values (new) - taken from what used to be row 4
$ws.Range("A2").Value = 111809606
$ws.Range("B2").Value = 56398
$ws.Range("E2").Value = 100109
$ws.Range("F2").Value = "Tretåig hackspett"
$ws.Range("G2").Value = "Picoides tridactylus"
$ws.Range("H2").Value = "(Linnaeus, 1758)"
$ws.Range("Q2").Value = 610564.5541715119
$ws.Range("R2").Value = 7180691.309759256
$ws.Range("Z2").Value = "14:33"
$ws.Range("AB2").Value = "14:33"

# Row 4 values (new) - taken from what used to be row 2
$ws.Range("A4").Value = 111809580
$ws.Range("B4").Value = 77515
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("Q4").Value = 610571.4165256479
$ws.Range("R4").Value = 7180702.680798599
$ws.Range("Z4").Value = "14:31"
$ws.Range("AB4").Value = "14:31"
